$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2..17, columns A..T (20 columns).
# Row r=2..13 are updates to existing rows (values recomputed after adding a new
# "Resolving-Mac" sending cluster); rows 14..17 are brand-new rows for the
# "Resolving-Mac" sending cluster.
$data = @(
    @("ECs","Fgf1","Fgfr1","ECs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,7.955277333333332,23.865832,0.05015625076675284,0.05015625076675283,9.471876707678222,85.246890369104,0.005700459244513984,0.005700459244513984),
    @("ECs","Fgf1","Fgfr1","FAPs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,82.48060333333333,247.44181,0.520021823355633,0.520021823355633,98.20476053986889,883.8428448588201,0.0591025677752937,0.0591025677752937),
    @("ECs","Fgf1","Fgfr1","MuSCs",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,66.90297433333333,200.708923,0.4218083439585467,0.4218083439585465,79.65740196222289,716.9166176600061,0.04794021157828463,0.04794021157828462),
    @("ECs","Fgf1","Fgfr1","Resolving-Mac",3,1,1.190640666666667,3.571922,0.1136540143525372,0.1136540143525372,3,1,1.271033333333333,3.8131,0.008013581919067616,0.008013581919067614,1.513343975355556,13.6200957782,0.0009107757544449438,0.0009107757544449436),
    @("FAPs","Fgf1","Fgfr1","ECs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,7.955277333333332,23.865832,0.05015625076675284,0.05015625076675283,41.11425482516356,370.028293426472,0.02474379061642757,0.02474379061642757),
    @("FAPs","Fgf1","Fgfr1","FAPs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,82.48060333333333,247.44181,0.520021823355633,0.520021823355633,426.2740821581123,3836.46673942301,0.256544516712841,0.256544516712841),
    @("FAPs","Fgf1","Fgfr1","MuSCs",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,66.90297433333333,200.708923,0.4218083439585467,0.4218083439585465,345.7661901712092,3111.895711540883,0.2080924547512396,0.2080924547512395),
    @("FAPs","Fgf1","Fgfr1","Resolving-Mac",3,1,5.168173666666667,15.504521,0.4933341355895272,0.4933341355895272,3,1,1.271033333333333,3.8131,0.008013581919067616,0.008013581919067614,6.56892100278889,59.12028902510001,0.003953373509019087,0.003953373509019086),
    @("MuSCs","Fgf1","Fgfr1","ECs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,7.955277333333332,23.865832,0.05015625076675284,0.05015625076675283,32.11699791846933,289.052981266224,0.01932897179098233,0.01932897179098233),
    @("MuSCs","Fgf1","Fgfr1","FAPs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,82.48060333333333,247.44181,0.520021823355633,0.520021823355633,332.9901968937133,2996.91177204342,0.2004034791412095,0.2004034791412095),
    @("MuSCs","Fgf1","Fgfr1","MuSCs",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,66.90297433333333,200.708923,0.4218083439585467,0.4218083439585465,270.1002865606874,2430.902579046186,0.1625544464934407,0.1625544464934407),
    @("MuSCs","Fgf1","Fgfr1","Resolving-Mac",3,1,4.037194,12.111582,0.3853751326204581,0.3853751326204581,3,1,1.271033333333333,3.8131,0.008013581919067616,0.008013581919067614,5.131408147133334,46.1826733242,0.003088235194825587,0.003088235194825587),
    @("Resolving-Mac","Fgf1","Fgfr1","ECs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,7.955277333333332,23.865832,0.05015625076675284,0.05015625076675283,0.6364407489804443,5.727966740823999,0.0003830291148289542,0.0003830291148289542),
    @("Resolving-Mac","Fgf1","Fgfr1","FAPs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,82.48060333333333,247.44181,0.520021823355633,0.520021823355633,6.598640721407778,59.38776649267,0.003971259726288792,0.003971259726288792),
    @("Resolving-Mac","Fgf1","Fgfr1","MuSCs",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,66.90297433333333,200.708923,0.4218083439585467,0.4218083439585465,5.352394053606777,48.171546482461,0.003221231135581728,0.003221231135581728),
    @("Resolving-Mac","Fgf1","Fgfr1","Resolving-Mac",1,0.3333333333333333,0.08000233333333333,0.240007,0.007636717437477471,0.007636717437477472,3,1,1.271033333333333,3.8131,0.008013581919067616,0.008013581919067614,0.1016856324111111,0.9151706917000001,0.00006119746077799785,0.00006119746077799783)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}

$wb.Save()
